$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.977.19"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "3.559.25"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "3.559.21"
$ws.Range("E7").Value = "  +1.53%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.12"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.429"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "4.157.13"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "3.550.65"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "67.767.36"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.46%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "449.34"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.627"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.67"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000133"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.52%  "
$ws.Range("D26").Value = "3.696.10"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.64"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.53"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.61"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.24%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "3.545.86"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.87"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.09"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "177.95"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.42%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0898"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.46"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.889"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.67"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.66"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.23%  "
